# Update TPM (transcripts per million)-derived NATMI metrics for Fbn1-Itgb1 ligand-receptor pairs.
# Ligand-expressing-cell counts moved from 2 -> 3 (detection rate 0.667 -> 1), which cascades through
# the average/total expression, specificity, and edge-weight columns recomputed with the new TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.445154666666667
$ws.Range("H2").Value = 10.335464
$ws.Range("I2").Value = 0.01110365039942287
$ws.Range("J2").Value = 0.01110365039942286
$ws.Range("M2").Value = 61.04160633333334
$ws.Range("N2").Value = 183.124819
$ws.Range("O2").Value = 0.2043613460574534
$ws.Range("P2").Value = 0.2043613460574534
$ws.Range("Q2").Value = 210.2977749201129
$ws.Range("R2").Value = 1892.679974281016
$ws.Range("S2").Value = 0.002269156941777437
$ws.Range("T2").Value = 0.002269156941777437

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.445154666666667
$ws.Range("H3").Value = 10.335464
$ws.Range("I3").Value = 0.01110365039942287
$ws.Range("J3").Value = 0.01110365039942286
$ws.Range("O3").Value = 0.3559304658284363
$ws.Range("P3").Value = 0.3559304658284363
$ws.Range("Q3").Value = 366.2697786740746
$ws.Range("R3").Value = 3296.428008066672
$ws.Range("S3").Value = 0.003952127459062683
$ws.Range("T3").Value = 0.003952127459062683

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.445154666666667
$ws.Range("H4").Value = 10.335464
$ws.Range("I4").Value = 0.01110365039942287
$ws.Range("J4").Value = 0.01110365039942286
$ws.Range("M4").Value = 131.3384093333333
$ws.Range("N4").Value = 394.015228
$ws.Range("O4").Value = 0.4397081881141102
$ws.Range("P4").Value = 0.4397081881141103
$ws.Range("Q4").Value = 452.4811338273102
$ws.Range("R4").Value = 4072.330204445792
$ws.Range("S4").Value = 0.004882365998582745
$ws.Range("T4").Value = 0.004882365998582745

# Row 5
$ws.Range("I5").Value = 0.938949437922138
$ws.Range("J5").Value = 0.938949437922138
$ws.Range("M5").Value = 61.04160633333334
$ws.Range("N5").Value = 183.124819
$ws.Range("O5").Value = 0.2043613460574534
$ws.Range("P5").Value = 0.2043613460574534
$ws.Range("Q5").Value = 17783.24879246736
$ws.Range("R5").Value = 160049.2391322062
$ws.Range("S5").Value = 0.1918849710136574
$ws.Range("T5").Value = 0.1918849710136574

# Row 6
$ws.Range("I6").Value = 0.938949437922138
$ws.Range("J6").Value = 0.938949437922138
$ws.Range("O6").Value = 0.3559304658284363
$ws.Range("P6").Value = 0.3559304658284363
$ws.Range("S6").Value = 0.334200710828975
$ws.Range("T6").Value = 0.3342007108289751

# Row 7
$ws.Range("I7").Value = 0.938949437922138
$ws.Range("J7").Value = 0.938949437922138
$ws.Range("M7").Value = 131.3384093333333
$ws.Range("N7").Value = 394.015228
$ws.Range("O7").Value = 0.4397081881141102
$ws.Range("P7").Value = 0.4397081881141103
$ws.Range("Q7").Value = 38262.81366888203
$ws.Range("R7").Value = 344365.3230199383
$ws.Range("S7").Value = 0.4128637560795055
$ws.Range("T7").Value = 0.4128637560795056

# Row 8
$ws.Range("G8").Value = 15.497141
$ws.Range("H8").Value = 46.491423
$ws.Range("I8").Value = 0.04994691167843914
$ws.Range("J8").Value = 0.04994691167843914
$ws.Range("M8").Value = 61.04160633333334
$ws.Range("N8").Value = 183.124819
$ws.Range("O8").Value = 0.2043613460574534
$ws.Range("P8").Value = 0.2043613460574534
$ws.Range("Q8").Value = 945.9703802141597
$ws.Range("R8").Value = 8513.733421927438
$ws.Range("S8").Value = 0.01020721810201856
$ws.Range("T8").Value = 0.01020721810201856

# Row 9
$ws.Range("G9").Value = 15.497141
$ws.Range("H9").Value = 46.491423
$ws.Range("I9").Value = 0.04994691167843914
$ws.Range("J9").Value = 0.04994691167843914
$ws.Range("O9").Value = 0.3559304658284363
$ws.Range("P9").Value = 0.3559304658284363
$ws.Range("Q9").Value = 1647.570269941706
$ws.Range("R9").Value = 14828.13242947535
$ws.Range("S9").Value = 0.01777762754039861
$ws.Range("T9").Value = 0.01777762754039861

# Row 10
$ws.Range("G10").Value = 15.497141
$ws.Range("H10").Value = 46.491423
$ws.Range("I10").Value = 0.04994691167843914
$ws.Range("J10").Value = 0.04994691167843914
$ws.Range("M10").Value = 131.3384093333333
$ws.Range("N10").Value = 394.015228
$ws.Range("O10").Value = 0.4397081881141102
$ws.Range("P10").Value = 0.4397081881141103
$ws.Range("Q10").Value = 2035.369848154382
$ws.Range("R10").Value = 18318.32863338944
$ws.Range("S10").Value = 0.02196206603602197
$ws.Range("T10").Value = 0.02196206603602197
